# Update the project journal worksheet:
#  - Row 12 (Monday 18 April 2022 entry): the logged time range, hours, and
#    details text are updated to reflect the actual session.
#  - A new row 13 is appended for the following session, "Tuesday 19 April
#    2022".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 12 ---
$ws.Range("C12").Value = "13:00-24:00"
$ws.Range("D12").Value = 11
$ws.Range("F12").Value = "Options functionality and persistant data."

# --- Add new row 13 ---
$ws.Range("B13").Value = "Tuesday 19 April 2022"
$ws.Range("C13").Value = "03:30-05:00"
$ws.Range("D13").Value = 1.5
$ws.Range("E13").Value = "Programming"
$ws.Range("F13").Value = "Save and high score system."

# --- Move the active selection to the new row, as the author did after typing ---
$ws.Range("D14").Select()
